$wb = $excel.ActiveWorkbook

# =====================================================================
# 1. "总计" (summary) sheet: insert a new row for 2022-Q4 at the top of
#    the data (row 2), pushing the existing Q3/Q2 rows down by one.
# =====================================================================
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
# The inserted row copies formatting down from the header row for some
# cells - strip that back to the default (unstyled) look used by the
# other data rows before writing values.
$summary.Range("A2:D2").ClearFormats()

# Row 2: 2022-Q4 (new)
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 10
$summary.Range("D2").Value = 4.14

# Row 3: 2022-Q3 (previously row 2)
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 18
$summary.Range("D3").Value = 5.77

# Row 4: 2022-Q2 (previously row 3)
$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 38
$summary.Range("D4").Value = 11.2

# Column A carries the same bordered/bold/centered style as the header
# row in this sheet - reapply it to the (re)written A2:A4 cells by
# copying the format from A3, which already has it.
$summary2 = $wb.Worksheets.Item("总计")
$summary2.Range("A3").Copy()
$summary3 = $wb.Worksheets.Item("总计")
$summary3.Range("A2").PasteSpecial(-4122)
$summary4 = $wb.Worksheets.Item("总计")
$summary4.Range("A3").Copy()
$summary5 = $wb.Worksheets.Item("总计")
$summary5.Range("A4").PasteSpecial(-4122)

# =====================================================================
# 2. Add the new "2022-Q4" sheet, positioned right after "总计" (i.e.
#    before "2022-Q3").
# =====================================================================
$q3ForAdd = $wb.Worksheets.Item("2022-Q3")
$q4Sheet = $wb.Worksheets.Add($q3ForAdd)
$q4Sheet.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

$rows = @(
    @(0, "010027", "景顺长城核心中景一年持有期混合", "44.17", "89.32", "4.00", "1.7668", 4),
    @(1, "010108", "景顺长城核心招景混合A",         "48.01", "89.61", "3.68", "1.7668", 7),
    @(2, "009190", "景顺长城核心优选一年持有期混合", "10.80", "90.64", "4.88", "0.5270", 5),
    @(3, "010783", "德邦沪港深龙头混合A",           "0.70",  "82.71", "3.58", "0.0251", 7),
    @(4, "013897", "德邦港股通成长精选混合型证券投资基金A", "0.46", "83.95", "3.61", "0.0166", 7),
    @(5, "013898", "德邦港股通成长精选混合型证券投资基金C", "0.42", "83.95", "3.61", "0.0152", 7),
    @(6, "010784", "德邦沪港深龙头混合C",           "0.35",  "82.71", "3.58", "0.0125", 7),
    @(7, "501303", "广发恒生中型股指数（LOF）A",     "0.24",  "90.77", "1.65", "0.0040", 4),
    @(8, "004996", "广发恒生中型股指数（LOF）C",     "0.12",  "90.77", "1.65", "0.0020", 4),
    @(9, "015752", "景顺长城核心招景混合C",         "0.01",  "89.61", "3.68", "0.0004", 7)
)

# --- header row (B1:H1), all text -----------------------------------
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 2   # B=2 .. H=8
    $cell = $q4Sheet.Cells.Item(1, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$i]
    $cell.ClearFormats()
}

# --- data rows 2..11 --------------------------------------------------
$r = 2
foreach ($row in $rows) {
    # Column A: numeric index (0-based)
    $q4Sheet.Cells.Item($r, 1).Value = $row[0]

    # Columns B-G: text (fund code / name / scale / position figures
    # stored verbatim as strings in the source data, incl. leading
    # zeros and fixed decimal formatting), force text interpretation.
    for ($c = 2; $c -le 7; $c++) {
        $cell = $q4Sheet.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c - 1]
        $cell.ClearFormats()
    }

    # Column H: numeric rank
    $q4Sheet.Cells.Item($r, 8).Value = $row[7]

    $r++
}

# --- styling: header row + column A use the bold/bordered/centered
#     "s=2" look already present elsewhere in the workbook; copy it
#     in from the equivalent cells on the neighbouring "2022-Q3" sheet.
$q3Style = $wb.Worksheets.Item("2022-Q3")
$q4Style = $wb.Worksheets.Item("2022-Q4")
$q3Style.Range("B1:H1").Copy()
$q4Style.Range("B1:H1").PasteSpecial(-4122)

$q3Style2 = $wb.Worksheets.Item("2022-Q3")
$q4Style2 = $wb.Worksheets.Item("2022-Q4")
$q3Style2.Range("A2").Copy()
$q4Style2.Range("A2:A11").PasteSpecial(-4122)
